$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the specific cell values per the diff
$ws.Range("A10").Value = -20.48049999999997
$ws.Range("A12").Value = -22.58920000000004
$ws.Range("B13").Value = 5.7801
$ws.Range("A18").Value = -22.54670000000004
